$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.343.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.072.31"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.18%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.54%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.067.26"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.522"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.153"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.480"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000247"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.98%  "
$ws.Range("E15").Value = "  -0.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.594.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.561.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.082.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +16.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "466.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.710"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0000102"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.06"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.114"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.09%  "
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.86"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "46.28"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.315"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.122"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0360"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.58%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "382.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.763.39"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.39%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.79%  "
